# Auto-generated script applying Bahamut_Profits leve-profit data refresh
# across all profession sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 761.5
$ws.Range("I18").Value = 837.8
$ws.Range("J18").Value = 380
$ws.Range("K18").Value = 837.8
$ws.Range("L18").Value = 380
$ws.Range("M18").Value = -553.8
$ws.Range("N18").Value = -948
$ws.Range("H19").Value = 1108.2258
$ws.Range("I19").Value = 1295.1177
$ws.Range("J19").Value = 881.2857
$ws.Range("K19").Value = 1295.1177
$ws.Range("L19").Value = 881.2857
$ws.Range("M19").Value = -1120.1177
$ws.Range("N19").Value = -1231.2857
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").Value = $null
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").Value = $null
$ws.Range("H29").Value = 2650
$ws.Range("I29").Value = 750
$ws.Range("J29").Value = 3192.8572
$ws.Range("K29").Value = 2250
$ws.Range("L29").Value = 9578.571599999999
$ws.Range("M29").Value = -1969
$ws.Range("N29").Value = -10140.5716
$ws.Range("H38").Value = 368.375
$ws.Range("I38").Value = 39.25
$ws.Range("J38").Value = 697.5
$ws.Range("K38").Value = 117.75
$ws.Range("L38").Value = 2092.5
$ws.Range("M38").Value = 254.25
$ws.Range("N38").Value = -2836.5
$ws.Range("H43").Value = 2132.1428
$ws.Range("I43").Value = 5633.3335
$ws.Range("J43").Value = 1177.2727
$ws.Range("K43").Value = 5633.3335
$ws.Range("L43").Value = 1177.2727
$ws.Range("M43").Value = -5564.3335
$ws.Range("N43").Value = -1315.2727
$ws.Range("H51").Value = 2879.6
$ws.Range("I51").Value = 1450
$ws.Range("J51").Value = 3832.6667
$ws.Range("K51").Value = 1450
$ws.Range("L51").Value = 3832.6667
$ws.Range("M51").Value = -966
$ws.Range("N51").Value = -4800.6667
$ws.Range("H58").Value = 1835.3636
$ws.Range("I58").Value = 1348.1666
$ws.Range("J58").Value = 2420
$ws.Range("K58").Value = 4044.4998
$ws.Range("L58").Value = 7260
$ws.Range("M58").Value = -3894.4998
$ws.Range("N58").Value = -7560
$ws.Range("H62").Value = 38876.93
$ws.Range("I62").Value = 74185.64
$ws.Range("J62").Value = 3568.2144
$ws.Range("K62").Value = 74185.64
$ws.Range("L62").Value = 3568.2144
$ws.Range("M62").Value = -73561.64
$ws.Range("N62").Value = -4816.2144
$ws.Range("H65").Value = 38876.93
$ws.Range("I65").Value = 74185.64
$ws.Range("J65").Value = 3568.2144
$ws.Range("K65").Value = 370928.2
$ws.Range("L65").Value = 17841.072
$ws.Range("M65").Value = -367808.2
$ws.Range("N65").Value = -24081.072
$ws.Range("H107").Value = 424.83334
$ws.Range("I107").Value = 431.46155
$ws.Range("K107").Value = 431.46155
$ws.Range("M107").Value = 1488.53845
$ws.Range("H138").Value = 2954.41
$ws.Range("I138").Value = 1009.3333
$ws.Range("J138").Value = 3673.822
$ws.Range("K138").Value = 3027.9999
$ws.Range("L138").Value = 11021.466
$ws.Range("M138").Value = 2112.0001
$ws.Range("N138").Value = -21301.466

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1574.7632
$ws.Range("I2").Value = 1335.3334
$ws.Range("J2").Value = 1985.2142
$ws.Range("K2").Value = 1335.3334
$ws.Range("L2").Value = 1985.2142
$ws.Range("M2").Value = -1222.3334
$ws.Range("N2").Value = -2211.2142
$ws.Range("H32").Value = 7366.013
$ws.Range("I32").Value = 6664.2266
$ws.Range("J32").Value = 60000
$ws.Range("K32").Value = 6664.2266
$ws.Range("L32").Value = 60000
$ws.Range("M32").Value = -6377.2266
$ws.Range("N32").Value = -60574
$ws.Range("H110").Value = 864.2105
$ws.Range("I110").Value = 782.2857
$ws.Range("J110").Value = 1093.6
$ws.Range("K110").Value = 782.2857
$ws.Range("L110").Value = 1093.6
$ws.Range("M110").Value = 1262.7143
$ws.Range("N110").Value = -5183.6
$ws.Range("H116").Value = 1574.7632
$ws.Range("I116").Value = 1335.3334
$ws.Range("J116").Value = 1985.2142
$ws.Range("K116").Value = 1335.3334
$ws.Range("L116").Value = 1985.2142
$ws.Range("M116").Value = 958.6666
$ws.Range("N116").Value = -6573.2142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1574.7632
$ws.Range("I3").Value = 1335.3334
$ws.Range("J3").Value = 1985.2142
$ws.Range("K3").Value = 1335.3334
$ws.Range("L3").Value = 1985.2142
$ws.Range("M3").Value = -1221.3334
$ws.Range("N3").Value = -2213.2142
$ws.Range("H107").Value = 5958.9585
$ws.Range("I107").Value = 370.9375
$ws.Range("J107").Value = 17135
$ws.Range("K107").Value = 370.9375
$ws.Range("L107").Value = 17135
$ws.Range("M107").Value = 1549.0625
$ws.Range("N107").Value = -20975

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2133.7144
$ws.Range("I31").Value = 1699.4073
$ws.Range("K31").Value = 1699.4073
$ws.Range("M31").Value = -1404.4073
$ws.Range("H34").Value = 2133.7144
$ws.Range("I34").Value = 1699.4073
$ws.Range("K34").Value = 1699.4073
$ws.Range("M34").Value = -1497.4073
$ws.Range("H41").Value = 13999.25
$ws.Range("I41").Value = 5000
$ws.Range("K41").Value = 5000
$ws.Range("M41").Value = -4572
$ws.Range("H50").Value = 9193.200000000001
$ws.Range("J50").Value = 9193.200000000001
$ws.Range("L50").Value = 9193.200000000001
$ws.Range("N50").Value = -10443.2
$ws.Range("H51").Value = 8510.375
$ws.Range("I51").Value = 7426.6665
$ws.Range("J51").Value = 9160.6
$ws.Range("K51").Value = 7426.6665
$ws.Range("L51").Value = 9160.6
$ws.Range("M51").Value = -6690.6665
$ws.Range("N51").Value = -10632.6
$ws.Range("H59").Value = 16281.2
$ws.Range("J59").Value = 16281.2
$ws.Range("L59").Value = 16281.2
$ws.Range("N59").Value = -18571.2
$ws.Range("H60").Value = 4310.5386
$ws.Range("I60").Value = 1000
$ws.Range("J60").Value = 8172.8335
$ws.Range("K60").Value = 1000
$ws.Range("L60").Value = 8172.8335
$ws.Range("M60").Value = -489
$ws.Range("N60").Value = -9194.833500000001
$ws.Range("H61").Value = 8510.375
$ws.Range("I61").Value = 7426.6665
$ws.Range("J61").Value = 9160.6
$ws.Range("K61").Value = 7426.6665
$ws.Range("L61").Value = 9160.6
$ws.Range("M61").Value = -7078.6665
$ws.Range("N61").Value = -9856.6
$ws.Range("H68").Value = 16413.572
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 16413.572
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 16413.572
$ws.Range("M68").Value = $null
$ws.Range("N68").Value = -17911.572
$ws.Range("H71").Value = 16413.572
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 16413.572
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 49240.716
$ws.Range("M71").Value = $null
$ws.Range("N71").Value = -56728.716
$ws.Range("H74").Value = 12516
$ws.Range("I74").Value = 5500
$ws.Range("J74").Value = 13685.333
$ws.Range("K74").Value = 5500
$ws.Range("L74").Value = 13685.333
$ws.Range("M74").Value = -4626
$ws.Range("N74").Value = -15433.333
$ws.Range("H77").Value = 12516
$ws.Range("I77").Value = 5500
$ws.Range("J77").Value = 13685.333
$ws.Range("K77").Value = 16500
$ws.Range("L77").Value = 41055.999
$ws.Range("M77").Value = -12132
$ws.Range("N77").Value = -49791.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 925
$ws.Range("I132").Value = 823.5294
$ws.Range("K132").Value = 7411.7646
$ws.Range("M132").Value = -4881.7646

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1092.88
$ws.Range("I97").Value = 1159.1305
$ws.Range("J97").Value = 331
$ws.Range("K97").Value = 1159.1305
$ws.Range("L97").Value = 331
$ws.Range("M97").Value = -663.1305
$ws.Range("N97").Value = -1323
$ws.Range("H107").Value = 344.18182
$ws.Range("I107").Value = 344.18182
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 344.18182
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1575.81818
$ws.Range("N107").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1884.5667
$ws.Range("I7").Value = 1844.8695
$ws.Range("J7").Value = 2015
$ws.Range("K7").Value = 1844.8695
$ws.Range("L7").Value = 2015
$ws.Range("M7").Value = -1732.8695
$ws.Range("N7").Value = -2239
$ws.Range("H124").Value = 52000
$ws.Range("J124").Value = 52000
$ws.Range("L124").Value = 52000
$ws.Range("N124").Value = -61820
$ws.Range("H126").Value = 1884.5667
$ws.Range("I126").Value = 1844.8695
$ws.Range("J126").Value = 2015
$ws.Range("K126").Value = 5534.6085
$ws.Range("L126").Value = 6045
$ws.Range("M126").Value = -3064.6085
$ws.Range("N126").Value = -10985
$ws.Range("H136").Value = 6815.4
$ws.Range("I136").Value = 3394.25
$ws.Range("K136").Value = 10182.75
$ws.Range("M136").Value = -7632.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1518.9565
$ws.Range("J107").Value = 1757.6666
$ws.Range("L107").Value = 5272.9998
$ws.Range("N107").Value = -9112.9998
$ws.Range("H109").Value = 14033.333
$ws.Range("J109").Value = 14033.333
$ws.Range("L109").Value = 14033.333
$ws.Range("N109").Value = -16807.333
$ws.Range("H132").Value = 1183.6451
$ws.Range("I132").Value = 962.7406999999999
$ws.Range("J132").Value = 2674.75
$ws.Range("K132").Value = 2888.2221
$ws.Range("L132").Value = 8024.25
$ws.Range("M132").Value = -358.2221
$ws.Range("N132").Value = -13084.25
